$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A376").Formula = '="2021-03-06"'
$ws.Range("A376").Copy()
$ws.Range("A376").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("B376").Value = 2182
$ws.Range("C376").Value = 104
$ws.Range("D376").Value = 1318
$ws.Range("E376").Value = 760
$ws.Range("F376").Value = 2783
$ws.Range("G376").Value = 50
$ws.Range("H376").Value = 482
$ws.Range("I376").Value = 2251
$ws.Range("J376").Value = 4.77
$ws.Range("K376").Value = 60.4
$ws.Range("L376").Value = 34.83
$ws.Range("M376").Value = 1.8
$ws.Range("N376").Value = 17.32
$ws.Range("O376").Value = 80.88

$ws.Range("A377").Formula = '="2021-03-07"'
$ws.Range("A377").Copy()
$ws.Range("A377").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("B377").Value = 2276
$ws.Range("C377").Value = 105
$ws.Range("D377").Value = 1266
$ws.Range("E377").Value = 905
$ws.Range("F377").Value = 2782
$ws.Range("G377").Value = 49
$ws.Range("H377").Value = 475
$ws.Range("I377").Value = 2258
$ws.Range("J377").Value = 4.61
$ws.Range("K377").Value = 55.62
$ws.Range("L377").Value = 39.76
$ws.Range("M377").Value = 1.76
$ws.Range("N377").Value = 17.07
$ws.Range("O377").Value = 81.16

$ws.Range("A378").Formula = '="2021-03-08"'
$ws.Range("A378").Copy()
$ws.Range("A378").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("B378").Value = 2166
$ws.Range("C378").Value = 100
$ws.Range("D378").Value = 1189
$ws.Range("E378").Value = 877
$ws.Range("F378").Value = 2780
$ws.Range("G378").Value = 46
$ws.Range("H378").Value = 472
$ws.Range("I378").Value = 2262
$ws.Range("J378").Value = 4.62
$ws.Range("K378").Value = 54.89
$ws.Range("L378").Value = 40.49
$ws.Range("M378").Value = 1.65
$ws.Range("N378").Value = 16.98
$ws.Range("O378").Value = 81.37

$ws.Range("A379").Formula = '="2021-03-09"'
$ws.Range("A379").Copy()
$ws.Range("A379").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("B379").Value = 2160
$ws.Range("C379").Value = 100
$ws.Range("D379").Value = 1315
$ws.Range("E379").Value = 745
$ws.Range("F379").Value = 2781
$ws.Range("G379").Value = 44
$ws.Range("H379").Value = 502
$ws.Range("I379").Value = 2235
$ws.Range("J379").Value = 4.63
$ws.Range("K379").Value = 60.88
$ws.Range("L379").Value = 34.49
$ws.Range("M379").Value = 1.58
$ws.Range("N379").Value = 18.05
$ws.Range("O379").Value = 80.37

$ws.Range("A380").Formula = '="2021-03-10"'
$ws.Range("A380").Copy()
$ws.Range("A380").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("B380").Value = 2176
$ws.Range("C380").Value = 98
$ws.Range("D380").Value = 1370
$ws.Range("E380").Value = 708
$ws.Range("F380").Value = 2780
$ws.Range("G380").Value = 43
$ws.Range("H380").Value = 500
$ws.Range("I380").Value = 2237
$ws.Range("J380").Value = 4.5
$ws.Range("K380").Value = 62.96
$ws.Range("L380").Value = 32.54
$ws.Range("M380").Value = 1.55
$ws.Range("N380").Value = 17.99
$ws.Range("O380").Value = 80.47

$ws.Range("A381").Formula = '="2021-03-11"'
$ws.Range("A381").Copy()
$ws.Range("A381").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("B381").Value = 2179
$ws.Range("C381").Value = 104
$ws.Range("D381").Value = 1345
$ws.Range("E381").Value = 730
$ws.Range("F381").Value = 2781
$ws.Range("G381").Value = 45
$ws.Range("H381").Value = 490
$ws.Range("I381").Value = 2246
$ws.Range("J381").Value = 4.77
$ws.Range("K381").Value = 61.73
$ws.Range("L381").Value = 33.5
$ws.Range("M381").Value = 1.62
$ws.Range("N381").Value = 17.62
$ws.Range("O381").Value = 80.76000000000001
